$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 through 85 in column C need to be updated to the new fitness value 7293.
# (Rows 86 and below already hold 7293 and are left untouched.)
$ws.Range("C2:C85").Value = 7293
